$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 3 (the "R" row) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 196
$wsOff.Range("C3").Value = 120
$wsOff.Range("D3").Value = 38
$wsOff.Range("E3").Value = 14
$wsOff.Range("F3").Value = 3
$wsOff.Range("G3").Value = 5

# --- DEF sheet: update row 3 (the "R" row) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 212
$wsDef.Range("C3").Value = 161
$wsDef.Range("D3").Value = 53
$wsDef.Range("E3").Value = 29
